# added 4wk low sales check
$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Forecast Comparison" ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# Row 2
$ws1.Range("H2").Value = 9.529999999999999
$ws1.Range("L2").Value = 0.93

# Row 3
$ws1.Range("D3").Value = 4
$ws1.Range("H3").Value = 8.529999999999999
$ws1.Range("L3").Value = 1.13

# Row 4
$ws1.Range("H4").Value = 7.53
$ws1.Range("L4").Value = 1.13

# Row 5
$ws1.Range("H5").Value = 6.53
$ws1.Range("L5").Value = 1.17

# Row 6
$ws1.Range("H6").Value = 5.53
$ws1.Range("L6").Value = 1.09

# Row 7
$ws1.Range("H7").Value = 4.76
$ws1.Range("L7").Value = 1.01

# Row 8
$ws1.Range("H8").Value = 3.76
$ws1.Range("L8").Value = 1.19

# Row 9
$ws1.Range("H9").Value = 2.76
$ws1.Range("L9").Value = 1.09

# Row 10
$ws1.Range("H10").Value = 1.76
$ws1.Range("L10").Value = 0.88

# Row 11
$ws1.Range("H11").Value = 0.76
$ws1.Range("I11").Value = "Low"
$ws1.Range("L11").Value = 0.93

# Row 12
$ws1.Range("L12").Value = 0.89

# Row 13
$ws1.Range("L13").Value = 0.88

# Row 14
$ws1.Range("D14").Value = 3
$ws1.Range("L14").Value = 1.04

# Row 15
$ws1.Range("D15").Value = 3
$ws1.Range("L15").Value = 0.95

# Row 16
$ws1.Range("D16").Value = 3
$ws1.Range("L16").Value = 0.84

# Row 17
$ws1.Range("D17").Value = 3
$ws1.Range("L17").Value = 0.99

# --- Sheet 2: "Summary" ---
$ws2 = $wb.Worksheets.Item("Summary")

# These cells hold numeric-looking text (t="inlineStr" in the source); force
# text formatting before assigning so they don't get reinterpreted as numbers.
$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "66"

$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "34"
